# Daily scrape update - 2025-09-28 03:20:54 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns C, D, H to fit the new (longer) content.
# (The host's ColumnWidth<->stored-width conversion adds ~5/6 of a
# character as cell padding, so back it out here to land on whole
# target widths of 90 / 56 / 57 in the saved OOXML.)
$padding = 5.0 / 6.0
$ws.Columns.Item(3).ColumnWidth = 90 - $padding
$ws.Columns.Item(4).ColumnWidth = 56 - $padding
$ws.Columns.Item(8).ColumnWidth = 57 - $padding

# Data rows: Opportunity ID, Title, Country, Applicants, Duration, Organization
# (Link is derived from the opportunity id)
$rows = @(
    @{ Row=2;  Id="1328063"; Title="[Accelerate Serbia] Architectural Designer Intern"; Country="Belgrade, Serbia"; Applicants="0 applicants"; Duration="9 - 12 Weeks"; Org="Structura Concept" },
    @{ Row=3;  Id="1328032"; Title="IT Analyst"; Country="2620 Ramada, Portugal"; Applicants="2 applicants"; Duration="9 - 12 Weeks"; Org="Pegadamotriz" },
    @{ Row=4;  Id="1327495"; Title="Content Creator"; Country="Sheraton Al Matar, El Nozha, Cairo Governorate, Egypt"; Applicants="5 applicants"; Duration="3 - 6 Months"; Org="Skyline Egypt Tours" },
    @{ Row=5;  Id="1327336"; Title="Partnerships and Innovation Intern"; Country="Hyderabad, Telangana, India"; Applicants="0 applicants"; Duration="9 - 12 Weeks"; Org="Arunodhaya Trust" },
    @{ Row=6;  Id="1326357"; Title="Human resource strategies to attract, engage, and inspire through the lens of marketing"; Country="日本、京都府京都市"; Applicants="26 applicants"; Duration="9 - 12 Weeks"; Org="Shinko Automotive Co., Ltd." },
    @{ Row=7;  Id="1325464"; Title="Accelerate Romania|Account Manager for Foreign Markets"; Country="Bucharest, Romania"; Applicants="28 applicants"; Duration="9 - 12 Weeks"; Org="Azuvioo" },
    @{ Row=8;  Id="1317664"; Title="SALES ATTENDED"; Country="Denizli, Kumkısık, Denizli, Türkiye"; Applicants="51 applicants"; Duration="6 - 18 Months"; Org="COTTON CASTLE TEKSTİL SANAYİ VE TİCARET ANONİM ŞİRKETİ" },
    @{ Row=9;  Id="1315600"; Title="Digital Marketing Executive"; Country="Cairo, Cairo Governorate, Egypt"; Applicants="13 applicants"; Duration="9 - 12 Weeks"; Org="KHEBRAT MISR" },
    @{ Row=10; Id="1314934"; Title="Social Media Marketing Executive"; Country="Petaling Jaya, Selangor, Malaysia"; Applicants="123 applicants"; Duration="6 - 18 Months"; Org="iWisers SDN BHD" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    # Opportunity ID is textual (values like "1328063"); force text so it
    # isn't coerced to a number, then reset style so no quote-prefix flag lingers.
    $ws.Cells.Item($rowNum, 1).Value = "'$($r.Id)"
    $ws.Cells.Item($rowNum, 1).Style = "Normal"
    $ws.Cells.Item($rowNum, 2).Value = "https://aiesec.org/opportunity/global-talent/$($r.Id)"
    $ws.Cells.Item($rowNum, 3).Value = $r.Title
    $ws.Cells.Item($rowNum, 4).Value = $r.Country
    $ws.Cells.Item($rowNum, 5).Value = "No"
    $ws.Cells.Item($rowNum, 6).Value = $r.Applicants
    $ws.Cells.Item($rowNum, 7).Value = $r.Duration
    $ws.Cells.Item($rowNum, 8).Value = $r.Org
}
